# Atualização de bases das ligas, do dia: 17-02-2024 às 22:47
#
# The underlying source data re-ordered a handful of fixtures (same date,
# different match) which made several adjacent rows swap places, and a
# later odds refresh tweaked a few more cells in isolation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($row1, $row2) {
    $r1 = $ws.Range("B$row1`:AC$row1")
    $r2 = $ws.Range("B$row2`:AC$row2")
    $v1 = $r1.Value()
    $v2 = $r2.Value()
    $r1.Value = $v2
    $r2.Value = $v1
}

# Rows whose entire contents (everything but the running index in column A)
# were swapped with their neighbour.
Swap-Rows 26 27
Swap-Rows 45 46
Swap-Rows 75 76
Swap-Rows 130 131
Swap-Rows 190 191

# Isolated odds-column corrections (closing line movement) on rows that
# were not otherwise touched.
$ws.Range("R196").Value = 1.99
$ws.Range("S196").Value = 1.91

$ws.Range("R197").Value = 1.93
$ws.Range("S197").Value = 1.97
$ws.Range("U197").Value = 1.9
$ws.Range("V197").Value = 1.95

$ws.Range("R198").Value = 2.02
$ws.Range("S198").Value = 1.88

$ws.Range("R199").Value = 2.08
$ws.Range("S199").Value = 1.82
$ws.Range("U199").Value = 1.9
$ws.Range("V199").Value = 1.95

$ws.Range("U200").Value = 2.05
$ws.Range("V200").Value = 1.8
